$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data rows (2-12): Candidate IDs and regenerated credentials ---
# Row 2
$ws.Range("A2").Value = 'WYOMj984'
$ws.Range("B2").Value = 231006230
$ws.Range("C2").Value = 'emztodg95'
$ws.Range("D2").Value = 'r2t$KQ%4'
$ws.Range("E2").Value = 'MR'
$ws.Range("F2").Value = 'paMwBcqk'
$ws.Range("G2").Value = 'bKmk'
$ws.Range("H2").Value = 'Candidate'

# Row 3
$ws.Range("A3").Value = 'DqjQI924'
$ws.Range("B3").Value = 231006229
$ws.Range("C3").Value = 'lqutyix61'
$ws.Range("D3").Value = 'Th#V3%7z'
$ws.Range("E3").Value = 'MR'
$ws.Range("F3").Value = 'sYexdNzS'
$ws.Range("G3").Value = 'foOZ'
$ws.Range("H3").Value = 'Candidate'

# Row 4
$ws.Range("A4").Value = 'WqNlL167'
$ws.Range("B4").Value = 231006228
$ws.Range("C4").Value = 'rczdvrk20'
$ws.Range("D4").Value = 'Nd2!4Ym%'
$ws.Range("E4").Value = 'MR'
$ws.Range("F4").Value = 'UKAxbWAS'
$ws.Range("G4").Value = 'nHNN'
$ws.Range("H4").Value = 'Candidate'

# Row 5
$ws.Range("A5").Value = 'OzBUq813'
$ws.Range("B5").Value = 231006227
$ws.Range("C5").Value = 'uvycwxq21'
$ws.Range("D5").Value = 'k!7h8$WJ'
$ws.Range("E5").Value = 'MR'
$ws.Range("F5").Value = 'rOYqRUyH'
$ws.Range("G5").Value = 'AUYe'
$ws.Range("H5").Value = 'Candidate'

# Row 6
$ws.Range("A6").Value = 'kQCsv913'
$ws.Range("B6").Value = 231006226
$ws.Range("C6").Value = 'imlavuu71'
$ws.Range("D6").Value = 'm%Vx28#J'
$ws.Range("E6").Value = 'MR'
$ws.Range("F6").Value = 'ZyVimdzP'
$ws.Range("G6").Value = 'KcjA'
$ws.Range("H6").Value = 'Candidate'

# Row 7
$ws.Range("A7").Value = 'wUTCV605'
$ws.Range("B7").Value = 231006225
$ws.Range("C7").Value = 'fofmdow86'
$ws.Range("D7").Value = 'R%H7du9&'
$ws.Range("E7").Value = 'MR'
$ws.Range("F7").Value = 'aEhXKuxB'
$ws.Range("G7").Value = 'AokR'
$ws.Range("H7").Value = 'Candidate'

# Row 8
$ws.Range("A8").Value = 'QwrpS898'
$ws.Range("B8").Value = 231006223
$ws.Range("C8").Value = 'lcdfuky80'
$ws.Range("D8").Value = 'S6&5eZ#c'
$ws.Range("E8").Value = 'MR'
$ws.Range("F8").Value = 'QnJMtvnW'
$ws.Range("G8").Value = 'BZAg'
$ws.Range("H8").Value = 'Candidate'

# Row 9
$ws.Range("A9").Value = 'qqjAA671'
$ws.Range("B9").Value = 231006221
$ws.Range("C9").Value = 'wligems53'
$ws.Range("D9").Value = 'D&$mSz84'
$ws.Range("E9").Value = 'MR'
$ws.Range("F9").Value = 'RYclXvcK'
$ws.Range("G9").Value = 'QNfz'
$ws.Range("H9").Value = 'Candidate'

# Row 10
$ws.Range("A10").Value = 'RMfNC544'
$ws.Range("B10").Value = 231006220
$ws.Range("C10").Value = 'meibsyi12'
$ws.Range("D10").Value = 'J7p%N9$g'
$ws.Range("E10").Value = 'MR'
$ws.Range("F10").Value = 'uzqGWExd'
$ws.Range("G10").Value = 'fWMF'
$ws.Range("H10").Value = 'Candidate'

# Row 11
$ws.Range("A11").Value = 'fzGGt246'
$ws.Range("B11").Value = 231006219
$ws.Range("C11").Value = 'pqkazif30'
$ws.Range("D11").Value = 'A4&h%3eG'
$ws.Range("E11").Value = 'MR'
$ws.Range("F11").Value = 'naLWMFBq'
$ws.Range("G11").Value = 'ZchR'
$ws.Range("H11").Value = 'Candidate'

# Row 12
$ws.Range("A12").Value = 'dBWRX203'
$ws.Range("B12").Value = 231006218
$ws.Range("C12").Value = 'sodfvpm25'
$ws.Range("D12").Value = 'js!%86QA'
$ws.Range("E12").Value = 'MR'
$ws.Range("F12").Value = 'tHZZplYC'
$ws.Range("G12").Value = 'bKiC'
$ws.Range("H12").Value = 'Candidate'

# --- Add new row 13 with matching format (bordered, like the other data rows) ---
$ws.Range("A13").Value = 'VzfWr715'
$ws.Range("B13").Value = 231006217
$ws.Range("C13").Value = 'ibhdtxe48'
$ws.Range("D13").Value = 's#U5$3vS'
$ws.Range("E13").Value = 'MR'
$ws.Range("F13").Value = 'JphdKHHB'
$ws.Range("G13").Value = 'aEFl'
$ws.Range("H13").Value = 'Candidate'

$newRow = $ws.Range("A13:H13")
$newRow.Borders.LineStyle = 1

# --- Refresh the sheet's used-range dimension / selection to cover the new row ---
$ws.Range("A1:H13").Select()
